$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A13").Value = "BanDoiSong"
$ws.Range("B13").Value = 12345
$ws.Range("C13").Value = "super"

[void]$ws.Range("C14").Select()
